# "1st day to add cookie"
#
# 1. Bump the auto "today" date field cached on the slide master and every
#    slide layout from 2017/6/17 -> 2017/6/18 (PowerPoint re-baked the
#    cached display text for the datetimeFigureOut field when the file was
#    reopened the next day).
# 2. Append two new slides (Title and Content layout) at the end of the
#    deck:
#      - slide 4: python snippet that configures selenium webdriver
#        options to ignore certificate errors, and resets
#        verificationErrors.
#      - slide 5: a new "cookies" section title slide.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "2017/6/17") {
                $tr.Characters(1, 9).Text = "2017/6/18"
            }
        }
    }
}

# Master + every layout carry their own cached copy of the date field.
Update-DatePlaceholder $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $p.SlideMaster.CustomLayouts.Item($li).Shapes
}

# --- New slide 4: selenium options / verificationErrors snippet ---------
$slide4 = $p.Slides.Add(4, 2)

$body4 = $slide4.Shapes.Item(2).TextFrame.TextRange
$body4.Text = "self.options.add_experimental_option"
$run = $body4.InsertAfter("(")
$run = $run.InsertAfter("`"")
$run = $run.InsertAfter("excludeSwitches")
$run = $run.InsertAfter("`"")
$run = $run.InsertAfter(", [")
$run = $run.InsertAfter("`"ignore-certificate-errors")
$run = $run.InsertAfter("`"])`rself.verificationErrors")
$run = $run.InsertAfter(" =")
$run = $run.InsertAfter(" []")

# --- New slide 5: "cookies" section title --------------------------------
$slide5 = $p.Slides.Add(5, 2)
$slide5.Shapes.Item(1).TextFrame.TextRange.Text = "cookies"
